$p = $ppt.ActivePresentation
Write-Host "Slide count before: $($p.Slides.Count)"
$p.Slides.Item(2).Delete()
Write-Host "Slide count after: $($p.Slides.Count)"
